# Adds a new user story (subscription bot feature) and a second user story
# (gamble feature) to the "User Story Specs" sheet, in rows 11 and 12,
# attributed to Joel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11 : subscription user story -------------------------------------
$ws.Range("A11").Value = "I, MasterTwitchUser, is a Twitch user that has just subscribed to the channel."
$ws.Range("D11").Value = "Logs name and can store the subscriptions in a data file for mining purposes."
$ws.Range("E11").Value = "Won't print out a message to the views showing the user subscribed becauser that is already built into Twitch."
$ws.Range("B11").Value = "The bot recognized the subscription of a user."
$ws.Range("C11").Value = "Should loga a message out of the chats view."
$ws.Range("F11").Value = "Joel"

# ---- Row 12 : gambling user story -------------------------------------------
$ws.Range("A12").Value = "I, TheGreatGabby01, is a Twitch user that wants to gamble some of their coins for potentailly more coins. They type ""!gamble xxx"" in chat. xxx = heads/tails"
$ws.Range("B12").Value = "The bot deducts coins from the user and then flips a coin to see if they won or not."
$ws.Range("C12").Value = "Access to the the data files where we keep all the users' coins and points."
$ws.Range("D12").Value = "A flashy message that shows the user gambling in chat for others to see."
$ws.Range("E12").Value = "A graphical image that moves to show that the user won or not."
$ws.Range("F12").Value = "Joel"

# ---- Formatting -------------------------------------------------------------
$ws.Range("A11:E12").WrapText = $true
$ws.Range("A11:E12").VerticalAlignment = -4108   # xlCenter
$ws.Range("F11:F12").WrapText = $true

$ws.Rows.Item(11).RowHeight = 135.75
$ws.Rows.Item(12).RowHeight = 119.25

# ---- Selection (matches the saved cursor position in the workbook) --------
$ws.Range("J12").Select()
